# Daily attendance processing - 2025-12-27 09:55:12
# Applies the day's attendance-recording results to the session-analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a percentage-looking value while keeping it stored as TEXT
# (the source sheet keeps these as literal strings like "55.7%", not numbers).
function Set-PercentText {
    param($addr, $text)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
}

# ---------------------------------------------------------------------------
# 1) Top summary panel (K/L columns, rows 6-10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 177     # Recorded Sessions
$ws.Range("L7").Value = 3       # Missing Sessions
Set-PercentText "L9"  "55.7%"   # Coverage %
Set-PercentText "L10" "75.4%"   # Average Attendance %

# ---------------------------------------------------------------------------
# 2) "Recorded By" column: reorder "<email>, System" -> "System, <email>"
#    for every session that was touched by the nightly System reconciliation.
# ---------------------------------------------------------------------------
$systemRows = @(8,9,10,34,35,36,60,61,62,86,87,88,112,113,114,138,139,140,
                164,167,170,191,194,197,218,221,224,245,248,251,272,275,278,
                299,302,305)
foreach ($r in $systemRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------------
# 3) Group Statistics table (rows 15-20) - recorded/missing counts & rates
#    shift because session #16 (27/12/2025) now has recorded attendance.
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 14
$ws.Range("P15").Value = 1
Set-PercentText "R15" "53.8%"
Set-PercentText "S15" "83.2%"

$ws.Range("O16").Value = 15
$ws.Range("P16").Value = 0
Set-PercentText "R16" "57.7%"
Set-PercentText "S16" "81.0%"

$ws.Range("O17").Value = 15
$ws.Range("P17").Value = 0
Set-PercentText "R17" "57.7%"
Set-PercentText "S17" "71.3%"

$ws.Range("O18").Value = 15
$ws.Range("P18").Value = 0
Set-PercentText "R18" "57.7%"
Set-PercentText "S18" "75.6%"

$ws.Range("O19").Value = 15
$ws.Range("P19").Value = 0
Set-PercentText "R19" "57.7%"
Set-PercentText "S19" "75.6%"

$ws.Range("O20").Value = 14
$ws.Range("P20").Value = 1
Set-PercentText "R20" "53.8%"
Set-PercentText "S20" "79.8%"

# ---------------------------------------------------------------------------
# 4) Session-log rows for 27/12/2025 (session #16) across every B1 group:
#    they were "Not Recorded" (pink) and are now "Recorded" (green), with
#    the recorder's email and the actual attendance count filled in.
# ---------------------------------------------------------------------------
function Set-SessionRecorded {
    param($row, $students)

    $ws.Range("G" + $row).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $row).Value = $students
    $ws.Range("I" + $row).Value = "Recorded"

    # Pink ("Not Recorded") -> Green ("Recorded") row highlight.
    $rng = $ws.Range("A" + $row + ":I" + $row)
    $rng.Interior.Color = 9498256   # RGB(144,238,144) == 0x90EE90
    $rng.Font.Color = 0
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

Set-SessionRecorded 16  "22/26"
Set-SessionRecorded 42  "24/27"
Set-SessionRecorded 68  "21/26"
Set-SessionRecorded 94  "21/27"
Set-SessionRecorded 120 "22/30"
Set-SessionRecorded 146 "18/23"
